$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.02266265261650346,0.1265126718081007,0.1266161815444171,1.757255418350454,1.137845425656764,1.103468705906124,0.9647071166073715,0.1560161159846984,1.292274900074801,0.4261677922716416,1.489690269959485),
    @(0.02131592053074627,0.1253780723251268,0.1264944989023018,1.754827232407379,1.133310009983759,1.107107687836915,0.961188425605819,0.1566015990387477,1.17583634719972,0.4029538030273159,1.506163100992381),
    @(0.02048099662832925,0.124715457766051,0.1264692536645455,1.75439569070889,1.13134163428839,1.109926329903615,0.9596665047277213,0.1570343233331108,1.104646113086687,0.3888669242011957,1.516851226035978),
    @(0.02013875061139458,0.1244540393270412,0.1264714224757135,1.7544859056364,1.130744296534999,1.111221771696663,0.9592065998947987,0.1572290740125375,1.075712745039056,0.3831685050675517,1.521350995495613),
    @(0.02008179995893755,0.124411151862823,0.1264725354544716,1.754516946647229,1.130657463484326,1.111445743557383,0.9591399059438714,0.1572625241554562,1.070913072933507,0.3822248356261113,1.522106892554245),
    @(0.02047638909559879,0.1247118972952492,0.1264692324550367,1.754395830488946,1.13133274985033,1.109943206347779,0.9596596536591377,0.1570368752640086,1.104255593650436,0.3887899026009762,1.516911327351554),
    @(0.02219997069365576,0.1261144240386187,0.12656396818549,1.756198109581007,1.136111869142937,1.104602092758881,0.9633611619781988,0.1562027881143351,1.25206440053455,0.4181291440474979,1.495251022444766),
    @(0.02551593656438911,0.1291331042586279,0.1271416386052771,1.768155148164752,1.1519860607699,1.098769932841861,0.9757012341211464,0.1551485148860827,1.5443007327126,0.476979624815101,1.457327289102981),
    @(0.02791299793319979,0.1315125375551247,0.1278043434882328,1.782102521822196,1.16765033263799,1.097324064215513,0.9878882777695281,0.154728915672095,1.760449981348643,0.5210169643015519,1.432238906276865),
    @(0.02899495225061344,0.1326297023171108,0.1281574535743317,1.789574841388784,1.175653620789404,1.097284904028811,0.9941152880643784,0.1546152451012617,1.859094704876838,0.5412241779775115,1.421427855011565),
    @(0.02940343357194308,0.133057701241647,0.1282985802906644,1.792567002989472,1.178811060656074,1.097359180564268,0.9965718883977246,0.1545833125015292,1.896493977349166,0.5489010792164351,1.41742056051153),
    @(0.02931551475958827,0.1329653044314725,0.1282678567238307,1.79191535114326,1.178125402163801,1.097339218283707,0.9960384262984121,0.1545896953931774,1.888437404213732,0.5472466185129292,1.418279749457227),
    @(0.02902858299982825,0.1326648150337846,0.1281689156973371,1.789817748091153,1.175910841309189,1.097289228096287,0.9943154166577273,0.1546123952298757,1.862170676113067,0.5418552640756218,1.421096436082784),
    @(0.02885266831181355,0.1324814003137433,0.1281092762420108,1.788554087524815,1.17457088508354,1.097270216502665,0.9932728700061233,0.1546277469346862,1.846087333466755,0.5385561361128453,1.422833020434567),
    @(0.02784211801053971,0.1314402236786805,0.1277823048443274,1.781636913150749,1.167145006675483,1.097339083946878,0.987495102258066,0.1547378987394268,1.754009618721511,0.5196998667482973,1.432957550010386),
    @(0.02721999757883964,0.1308103641344545,0.1275949349344856,1.777682544369938,1.162814639980809,1.097539873295716,0.9841258212111015,0.1548252547314632,1.697603504639631,0.5081766748842043,1.439322806689333),
    @(0.02686137279351186,0.1304513588972682,0.1274920263935044,1.775514205323901,1.16040649023428,1.097713575553968,0.9822521723540447,0.1548827661023182,1.665190127671053,0.5015652682322411,1.443040571094826),
    @(0.02673981198005748,0.1303303694308369,0.1274580187722378,1.77479825491038,1.159605294686472,1.097782381683913,0.9816288154895076,0.1549034862048089,1.65422068526243,0.4993295876365949,1.444309067373347),
    @(0.02728630599111881,0.130877075302962,0.1276143777218905,1.778092508795538,1.16326706578478,1.097512472917529,0.9844778324865047,0.1548152034103651,1.703604941456717,0.5094016395780443,1.438639351949526),
    @(0.02911289532325156,0.1327529419833269,0.1281977760414748,1.790429449934251,1.176557866420751,1.097301491872827,0.9948188285406019,0.1546054260833216,1.869884646648416,0.543438162327945,1.420266755534946),
    @(0.03029949469727455,0.1340077782738405,0.128622254384716,1.799440049772087,1.185983346889827,1.097683053205515,1.002151949853342,0.1545330982544186,1.978817891913195,0.5658278504632221,1.408764107929649),
    @(0.02966684532729857,0.1333354232810322,0.1283917548139044,1.794544072401564,1.180884963750259,1.097431825098994,0.9981854276658595,0.1545657709121073,1.920654713407259,0.5538648678195557,1.414857060153878),
    @(0.02725633094537017,0.1308469054954102,0.1276055726451872,1.777906836581863,1.163062270449615,1.097524679136157,0.984318490735717,0.1548197249078385,1.700891644301521,0.508847791192153,1.438948160252416),
    @(0.02462574747838886,0.1282879497652658,0.1269434731496553,1.764015876577687,1.146991683314042,1.099849803755262,0.9718165438495205,0.1553714274959148,1.464988934601934,0.476979624815101,1.467099449726749)
)

$colLetters = @("C","D","E","F","G","H","I","J","K","M","N")
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $ws.Range($colLetters[$j] + $rowNum).Value = $rowVals[$j]
    }
}
